$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 203

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 98

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 88

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 56

$ws.Range("A6:B6").Delete()
